# Report regenerated for archive: localization status moved from
# "Ready for handoff" to "In Translation" on the Overview sheet and on
# each per-locale sheet's Status column. Re-running the report's
# column autosizing also narrows the affected Status columns to match
# the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn (E) and de-de (F) status columns, data rows 2-4
$overview.Range("E2:F4").Value = "In Translation"

# Per-locale sheets: Status column is column C, data rows 2-4
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter. The target
# autofit width (matching the regenerated report) is 13.4101845877511
# character-units; the host's ColumnWidth setter quantizes to whole
# pixels, so 12.5 is the input that lands closest to that target after
# quantization.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
